$d = $word.ActiveDocument

$map = @(
    @("525×9=4725", "489×8=3912"),
    @("613×8=4904", "430×7=3010"),
    @("340×2=680", "611×2=1222"),
    @("973×6=5838", "488×6=2928"),
    @("445×8=3560", "477×9=4293"),
    @("879×5=4395", "816×5=4080"),
    @("365×8=2920", "467×3=1401"),
    @("663×2=1326", "130×2=260"),
    @("530×8=4240", "234×9=2106"),
    @("228×7=1596", "325×3=975"),
    @("360×3=1080", "701×9=6309"),
    @("332×9=2988", "675×6=4050"),
    @("239×6=1434", "198×4=792"),
    @("669×8=5352", "942×3=2826"),
    @("535×4=2140", "811×8=6488"),
    @("197×2=394", "335×8=2680"),
    @("959×3=2877", "362×4=1448"),
    @("953×4=3812", "259×2=518"),
    @("238×8=1904", "665×6=3990"),
    @("839×2=1678", "616×7=4312"),
    @("106×4=424", "341×2=682"),
    @("739×8=5912", "723×9=6507"),
    @("112×9=1008", "693×6=4158"),
    @("768×8=6144", "655×9=5895"),
    @("128×2=256", "140×8=1120")
)

foreach ($pair in $map) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
